$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "232.38", "1.00") need to be
# forced to Text format first, otherwise Excel auto-converts the literal into a
# numeric value (dropping the significant trailing zeros the source data relies on).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated coin price / volume figures scraped this run.
$ws.Range("D2").Value = "42.122.01"
$ws.Range("E2").Value = "  -4.20%  "
$ws.Range("D3").Value = "2.237.56"
$ws.Range("E3").Value = "  -4.90%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "232.38"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  -6.28%  "
$ws.Range("D7").Value = "68.64"
$ws.Range("E7").Value = "  -5.41%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").Value = "  -6.20%  "
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "58.32"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "35.75"
$ws.Range("E12").Value = "  +8.89%  "
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").Value = "6.71"
$ws.Range("E14").Value = "  -7.43%  "
$ws.Range("D15").Value = "2.573.05"
$ws.Range("E15").Value = "  -4.80%  "
$ws.Range("D16").Value = "14.92"
$ws.Range("E16").Value = "  -8.89%  "
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  -5.64%  "
$ws.Range("D18").Value = "2.228.74"
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").Value = "42.054.35"
$ws.Range("E19").Value = "  -4.19%  "
$ws.Range("D20").Value = "0.0₃0967"
$ws.Range("E20").Value = "  -5.96%  "
$ws.Range("D21").Value = "73.14"
$ws.Range("E21").Value = "  -6.51%  "
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  -7.61%  "
$ws.Range("D23").Value = "234.92"
$ws.Range("E23").Value = "  -7.72%  "
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -6.31%  "
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  -4.93%  "
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("D30").Value = "168.97"
$ws.Range("E30").Value = "  -4.81%  "
$ws.Range("D31").Value = "20.48"
$ws.Range("E31").Value = "  -8.82%  "
$ws.Range("E32").Value = "  -7.30%  "
$ws.Range("D33").Value = "0.126"
$ws.Range("E33").Value = "  -7.34%  "
$ws.Range("D34").Value = "0.0710"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D35").Value = "5.23"
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("D36").Value = "4.68"
$ws.Range("E36").Value = "  -8.68%  "
$ws.Range("D37").Value = "3.59"
$ws.Range("E37").Value = "  -4.30%  "
$ws.Range("D38").Value = "21.70"
$ws.Range("E38").Value = "  +14.95%  "
$ws.Range("E39").Value = "  -5.89%  "
$ws.Range("D40").Value = "5.98"
$ws.Range("E40").Value = "  -7.03%  "
$ws.Range("D41").Value = "0.0264"
$ws.Range("E41").Value = "  -3.79%  "
$ws.Range("D42").Value = "65.60"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "4.88"
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("D44").Value = "8.88"
$ws.Range("E44").Value = "  -4.11%  "
$ws.Range("D45").Value = "0.0997"
$ws.Range("E45").Value = "  -8.27%  "
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.188"
$ws.Range("E47").Value = "  -5.80%  "
$ws.Range("D48").Value = "1.17"
$ws.Range("E48").Value = "  -5.27%  "
$ws.Range("D49").Value = "4.31"
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  -6.58%  "
$ws.Range("D51").Value = "9.96"
$ws.Range("E51").Value = "  +6.22%  "
